$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Trial 19", 0, 0, 0, 10.46670031547546, 10.3735408782959, 0, 10.28624629974365),
    @("Trial 20", 0, 0, 0, 10.47215461730957, 10.36481404304504, 0, 10.29764723777771),
    @("Trial 21", 0, 0, 0, 21.09299921989441, 20.84064340591431, 0, 20.98466849327087),
    @("Trial 22", 0, 0, 0, 21.25983667373657, 20.84869337081909, 0, 20.91546201705933),
    @("Trial 23", 0, 0, 0, 21.23253870010376, 20.91148519515991, 0, 20.95027732849121)
)

$startRow = 20
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

$wb.Save()
